# update database and change read_price algorithm
#
# The quarterly "dollar" income-statement sheet keeps one column per
# fiscal quarter (columns D..M). A new quarter's figures were published,
# so the oldest quarter (column D) is dropped, every remaining quarter
# shifts one column to the left, and the newly published quarter's data
# is appended as the new last column (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# 1) Drop the oldest quarter (column D) -- this shifts E:M left to D:L,
#    carrying along both the values and the per-column formatting.
$ws.Columns.Item(4).Delete()

# 2) Clone the formatting of the (now) last populated column (L) into
#    the new column (M) so the freshly-added quarter matches the
#    existing look (borders/fills/alignment/number formats).
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new quarter's header, publish date, and figures.
$ws.Range("M8").Value = "فصل دوم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-11 (2)"

$ws.Range("M11").Value = 2059
$ws.Range("M12").Value = -1042
$ws.Range("M13").Value = 1017
$ws.Range("M14").Value = -410
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 260
$ws.Range("M17").Value = 866
$ws.Range("M18").Value = -32
$ws.Range("M19").Value = -18
$ws.Range("M20").Value = 816
$ws.Range("M21").Value = -171
$ws.Range("M22").Value = 645
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 645
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 4840
$ws.Range("M27").Value = 0

# 4) A couple of figures in the existing columns were restated by the
#    updated read_price algorithm (not just shifted).
$ws.Range("K9").Value = "1402-02-11 (5)"
$ws.Range("K23").Value = 1
$ws.Range("K24").Value = 478
